$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.520797848701477
$ws.Range("B1").Value = 1.714658617973328
$ws.Range("C1").Value = 2.093204975128174
$ws.Range("D1").Value = 2.809376955032349
$ws.Range("E1").Value = 6.489599704742432
